# regen sval data to filter save games
# Update the per-game stat columns (TB, d2S, K, IP) and the derived "sum" column (G)
# for every data row (2-15) of the active sheet. Column F (Win) and column A (date)
# are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548),
    @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 3.56341032713086),
    @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447),
    @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387),
    @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 3.56341032713086),
    @(0.2917716402565462, 0.306821227259698,  3.537761648806719, 10.19245300693656,  14.32880752325952),
    @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(0.1190320826869504, 0.306821227259698,  3.537761648806719, 0.4942365360607697, 4.457851494814137),
    @(1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 7.143138311642302),
    @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548),
    @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 3.56341032713086),
    @(0.6606524410359556, 0.306821227259698,  0.1494219747398047, 0.4942365360607697, 1.611132179096228)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
